# Generate Report for Handoff
#
# Refreshes the "Latest Handoff Datetime" column (D) on the per-locale
# status sheets for the source files whose handoff was just (re-)generated.
# Rows 2, 3, 5, and 11 already carry their own correct/unique handoff
# timestamps and are left untouched; rows 4, 6, 7, 8, 9 and 10 previously
# shared a stale placeholder timestamp and now get the freshly generated
# handoff datetime for this run.

$wb = $excel.ActiveWorkbook

$rowsToUpdate = 4, 6, 7, 8, 9, 10

$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rowsToUpdate) {
    $zhcn.Cells.Item($r, 4).Value = "2016-03-03 11:15:28"
}

$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rowsToUpdate) {
    $dede.Cells.Item($r, 4).Value = "2016-03-03 11:15:42"
}
